# Apply the "Shot-Boundary-Detection ground truth" restructuring:
#  - insert 3 new leading columns (Match, Game half, Game time (min:sec))
#  - strip the "test_" artifact prefix from what is now the League column
#  - re-point the workbook metadata the same way Excel would after a resave

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert three blank columns before the current column A.
#    (old A:F -> new D:I, all widths/formats travel with the cells)
# ---------------------------------------------------------------------
$ws.Range("A1:C1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2. Strip the "test_" dataset-artifact prefix from the (now) League
#    column values, and rewrite its header from "video" to "League".
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "League"
$ws.Range("D2").Value = "bundesliga"
$ws.Range("D3").Value = "epl"
$ws.Range("D4").Value = "france"
$ws.Range("D5").Value = "italy"
$ws.Range("D6").Value = "UEFA"

# ---------------------------------------------------------------------
# 3. New header row cells (A1:C1) - use the same "bold header" look as
#    the other header cells, but with the thinner/plain font used for
#    this group of columns.
# ---------------------------------------------------------------------
$ws.Range("D1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)

$ws.Range("A1").Value = "Match"
$ws.Range("B1").Value = "Game half"
$ws.Range("C1").Value = "Game time (min:sec)"

$ws.Range("A1:C1").Font.Size = 11
$ws.Range("A1:C1").Font.Bold = $false
$ws.Range("A1:C1").Font.ThemeColor = 1
$ws.Range("A1:C1").HorizontalAlignment = -4108
$ws.Range("A1:C1").VerticalAlignment = -4108
$ws.Range("A1:C1").WrapText = $true

# ---------------------------------------------------------------------
# 4. New data cells (A2:C6) - same look as the existing data columns,
#    but sized/colored like the new header font family.
# ---------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("A2:C6").PasteSpecial(-4122)

$ws.Range("A2:C6").Font.Size = 14
$ws.Range("A2:C6").Font.ThemeColor = 1
$ws.Range("A2:C6").HorizontalAlignment = -4108
$ws.Range("A2:C6").VerticalAlignment = -4108
$ws.Range("A2:C6").WrapText = $true

$ws.Range("A2").Value = "2015-02-21 - 17-30 Paderborn 0 - 6 Bayern Munich"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("C2").Value = "17:55"

$ws.Range("A3").Value = "2015-02-21 - 18-00 Crystal Palace 1 - 2 Arsenal"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1"
$ws.Range("C3").Value = "45:08"

$ws.Range("A4").Value = "2015-04-05 - 22-00 Marseille 2 - 3 Paris SG"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2"
$ws.Range("C4").Value = "05:48"

$ws.Range("A5").Value = "2016-08-21 - 21-45 Pescara 2 - 2 Napoli"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "1"
$ws.Range("C5").Value = "34:54"

$ws.Range("A6").Value = "2017-03-08 - 22-45 Barcelona 6 - 1 Paris SG"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "1"
$ws.Range("C6").Value = "39:24"

# ---------------------------------------------------------------------
# 5. Header row / data rows for D:H also gained center/middle alignment
#    in the resave; line them up with the new look.
# ---------------------------------------------------------------------
$ws.Range("D1:H1").HorizontalAlignment = -4108
$ws.Range("D1:H1").VerticalAlignment = -4108
$ws.Range("D2:H6").HorizontalAlignment = -4108
$ws.Range("D2:H6").VerticalAlignment = -4108

$ws.Range("A1").Select()

# ---------------------------------------------------------------------
# 6. Workbook-level metadata touch-ups (mirrors what a newer Excel
#    build stamps on save / a change of authoring machine).
# ---------------------------------------------------------------------
$wb.Title = $wb.Title
